$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D31").Value = "[1, 0, 0, 0, 0, 0, 1]"
$ws.Range("E31").Value = "['Normal', 'SoftwareFault']"

$ws.Range("D36").Value = "[1, 0, 1, 0, 0, 0, 0]"
$ws.Range("E36").Value = "['Normal', 'HardwareFault']"

$ws.Range("D54").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E54").Value = "['Normal']"

$ws.Range("D65").Value = "[0, 0, 0, 0, 0, 0, 0]"
$ws.Range("E65").Value = "[]"

$ws.Range("D74").Value = "[1, 0, 1, 0, 0, 0, 1]"
$ws.Range("E74").Value = "['Normal', 'HardwareFault', 'SoftwareFault']"
